# Reordered tables in ER diagram
#
# Slide 3 holds the ER diagram with two small "entity" tables
# ("Table 5" = customers, "Table 7" = purchases) placed side by side.
# This swaps their horizontal placement (left <-> right), leaving their
# vertical position, size and all other formatting untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

$table5 = $s.Shapes.Item("Table 5")
$table7 = $s.Shapes.Item("Table 7")

# Target EMU offsets (converted to points, 1 pt = 12700 EMU) taken from
# the canonical edit: Table 5 moves from x=4916285 EMU to x=1108559 EMU,
# Table 7 moves from x=1159161 EMU to x=4890984 EMU. Vertical position is
# unchanged for both shapes.
$table5.Left = 87.28812
$table7.Left = 385.11685
